# Insert a new weekly data row at row 250 (pushing existing rows 250-374 down
# to 251-375), and populate it with the new observation.
#
# The new row carries the same Mercado/Region/Categoria/etc. values as the
# (now shifted) row below it, but with its own Fecha, Volumen, Precio
# minimo/maximo/promedio ponderado and Precio $/Kg values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 250:374 down to 251:375, creating a blank row 250.
$ws.Rows(250).Insert()

# Copy formatting (number format, style) from the row that is now directly
# below (old row 250, now row 251) onto the freshly inserted row, so the new
# row matches the rest of the table (e.g. the date-style cell in column D).
$ws.Range("A251:R251").Copy()
$ws.Range("A250:R250").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$newRow = 250

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 44917
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112037
$ws.Cells.Item($newRow, 7).Value = "Cebollín"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 70
$ws.Cells.Item($newRow, 11).Value = 7000
$ws.Cells.Item($newRow, 12).Value = 7000
$ws.Cells.Item($newRow, 13).Value = 7000
$ws.Cells.Item($newRow, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 194
$ws.Cells.Item($newRow, 17).Value = 36
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
